$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header of column C from "cuenta" to "cuenta prestamo"
$ws.Range("C1").Value = "cuenta prestamo"

# Column C auto-fits wider to accommodate the new, longer header text
$ws.Columns.Item(3).ColumnWidth = 14

# Move the active selection to D1 (as reflected in the saved file)
$ws.Range("D1").Select()
